$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$nl = [char]10

# --- Update existing rows 9-13 (column F / "GB") ---
$ws.Range("F9").Value  = "I2C / SPI"
$ws.Range("F10").Value = "Orientation de la voile"
$ws.Range("F11").Value = "Tension batterie"
$ws.Range("F12").Value = "Test Complet"
$ws.Range("F13").Value = "Girouette"

# --- Append new rows 15-25 ---

# Row 15
$ws.Range("B15").Value = "SingleUseId16"
$ws.Range("C15").Value = "Default"
$ws.Range("D15").Value = "Center"
$ws.Range("E15").Value = "LTR"
$ws.Range("F15").Value = "Utilisez l'analyseur logique" + $nl + "pour visualiser les communications" + $nl + "entre le micro et les peripheriques" + $nl + "I2C et SPI"

# Row 16
$ws.Range("B16").Value = "SingleUseId17"
$ws.Range("C16").Value = "Default"
$ws.Range("D16").Value = "Left"
$ws.Range("E16").Value = "LTR"
$ws.Range("F16").Value = "Orientation bateau: "

# Row 17
$ws.Range("B17").Value = "SingleUseId18"
$ws.Range("C17").Value = "Default"
$ws.Range("D17").Value = "Center"
$ws.Range("E17").Value = "LTR"
$ws.Range("F17").Value = "Servo"

# Row 18
$ws.Range("B18").Value = "SingleUseId19"
$ws.Range("C18").Value = "Default"
$ws.Range("D18").Value = "Center"
$ws.Range("E18").Value = "LTR"
$ws.Range("F18").Value = "<value>"

# Row 19
$ws.Range("B19").Value = "SingleUseId20"
$ws.Range("C19").Value = "Default"
$ws.Range("D19").Value = "Left"
$ws.Range("E19").Value = "LTR"
$ws.Range("F19").Value = "--" + [char]0xB0 + " (-- ms)"

# Row 20
$ws.Range("B20").Value = "SingleUseId21"
$ws.Range("C20").Value = "Large"
$ws.Range("D20").Value = "Center"
$ws.Range("E20").Value = "LTR"
$ws.Range("F20").Value = "<value>"

# Row 21
$ws.Range("B21").Value = "SingleUseId22"
$ws.Range("C21").Value = "Large"
$ws.Range("D21").Value = "Left"
$ws.Range("E21").Value = "LTR"
$ws.Range("F21").Value = "0.0V" + $nl + "[0x000]"

# Row 22
$ws.Range("B22").Value = "SingleUseId23"
$ws.Range("C22").Value = "Large"
$ws.Range("D22").Value = "Center"
$ws.Range("E22").Value = "LTR"
$ws.Range("F22").Value = "<value>"

# Row 23
$ws.Range("B23").Value = "SingleUseId24"
$ws.Range("C23").Value = "Large"
$ws.Range("D23").Value = "Left"
$ws.Range("E23").Value = "LTR"
# Force text storage so "0" isn't stored as a number (matches source data type)
$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = "0"

# Row 24
$ws.Range("B24").Value = "SingleUseId25"
$ws.Range("C24").Value = "Default"
$ws.Range("D24").Value = "Left"
$ws.Range("E24").Value = "LTR"
$ws.Range("F24").Value = "Communication"

# Row 25
$ws.Range("B25").Value = "SingleUseId26"
$ws.Range("C25").Value = "Default"
$ws.Range("D25").Value = "Center"
$ws.Range("E25").Value = "LTR"
$ws.Range("F25").Value = "Rotation" + $nl + "plateau"

Write-Output "edit applied"
